# Update logical model: StructureDefinition-log-cercle-de-soins
# (commit: update logical model 800080629e233131db6c94e47e63653bef40085d)
$wb = $excel.ActiveWorkbook

# ---- Metadata sheet: bump the generation Date ----
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = '2024-02-14T09:17:54+00:00'

# ---- Elements sheet: fix placeholder Short/Definition text, add new row ----
$ws = $wb.Worksheets.Item("Elements")

# idCercleSoins (row 3)
$ws.Range("L3").Value = 'Identifiant du cercle de soins.'
$ws.Range("M3").Value = 'Identifiant du cercle de soins.'

# dateCreation (row 4)
$ws.Range("L4").Value = 'Date de création du cercle de soin.'
$ws.Range("M4").Value = 'Date de création du cercle de soin.'

# dateMAJ (row 5) - Short becomes the new text; Definition keeps the long-standing note
# that used to live in the Short cell.
$ws.Range("M5").Value = 'Le concept de cercle de soins a plusieurs dates de mise à jour mais chaque version de la ressource, et donc chaque instance de la ressource, ne peut avoir qu’une seule date de mise à jour'
$ws.Range("L5").Value = 'Date de mise à jour du cercle de soin.'

# dateFin (row 6)
$ws.Range("L6").Value = 'Date de fin d''existence du cercle de soins.'
$ws.Range("M6").Value = 'Date de fin d''existence du cercle de soins.'

# statut (row 7)
$ws.Range("L7").Value = 'Statut du cercle de soins.'
$ws.Range("M7").Value = 'La liste des valeurs possibles n''est pas définie par ces spécifications. Les codes possibles ainsi que le sens donné sont définis par le gestionnaire en fonction du projet. Il peut s''inspirer du jeu de valeur FHIR CareTeamStatus (proposed | active | suspended | inactive | entered-in-error).]]'

# metadonnee (row 8)
$ws.Range("L8").Value = 'Informations relatives à la gestion des classes et des données.'
$ws.Range("M8").Value = 'Informations relatives à la gestion des classes et des données.'

# ---- New row 9: PersonnePriseCharge (clone formatting from row 8, then fill values) ----
$ws.Range("A8:AK8").Copy()
$ws.Range("A9:AK9").PasteSpecial(-4122)

$ws.Range("A9").Value = 'log-cercle-de-soins.PersonnePriseCharge'
$ws.Range("B9").Value = 'log-cercle-de-soins.PersonnePriseCharge'
$ws.Range("F9").Value = "1"
$ws.Range("G9").Value = "1"
$ws.Range("K9").Value = 'http://interop.esante.gouv.fr/ig/fhir/cds/StructureDefinition/log-personne-prise-charge
'
$ws.Range("L9").Value = 'Personne prise en charge.'
$ws.Range("M9").Value = 'La personne prise en charge.'
$ws.Range("AF9").Value = 'log-cercle-de-soins.PersonnePriseCharge'
$ws.Range("AG9").Value = "1"
$ws.Range("AH9").Value = "1"

